$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Payment" column header to "sbp_payment"
$ws.Range("B1").Value = "sbp_payment"

# Replace flat yearly totals with formulas summing the underlying
# soil/climate/census-derived payment records for each year
$ws.Range("B15").Formula = "=50.72+50.72+51.82"
$ws.Range("B16").Formula = "=65.86+66.24"
$ws.Range("B17").Formula = "=40+43.58+43.45"
$ws.Range("B18").Formula = "=69.18+68.78"

# Highlight the 2011 payment cell with a light accent fill
$ws.Range("B17").Interior.Color = 14461583

# No payment data available for 2013 and 2014
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0

# Restore the sheet selection to the newly edited range
$ws.Range("B19:B24").Select() | Out-Null
